$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.637.61"
$ws.Range("E2").Value = "  -4.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.770.41"
$ws.Range("E3").Value = "  -5.60%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.39"
$ws.Range("E5").Value = "  -5.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.44"
$ws.Range("E6").Value = "  +7.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("E7").Value = "  -5.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.734"
$ws.Range("E9").Value = "  -3.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.01"
$ws.Range("E11").Value = "  -7.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000310"
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.00"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.381.28"
$ws.Range("E14").Value = "  -5.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.829.51"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.89"
$ws.Range("E16").Value = "  -4.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.37"
$ws.Range("E17").Value = "  -6.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.17"
$ws.Range("E18").Value = "  -8.33%  "
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.561.71"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "419.93"
$ws.Range("E21").Value = "  -5.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  -6.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "91.15"
$ws.Range("E23").Value = "  -5.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.14"
$ws.Range("E24").Value = "  -7.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.28"
$ws.Range("E25").Value = "  -8.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.13"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.83"
$ws.Range("E27").Value = "  -8.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  -5.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.50"
$ws.Range("E30").Value = "  +8.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.86"
$ws.Range("E31").Value = "  -6.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.16"
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "46.03"
$ws.Range("E33").Value = "  -6.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.121"
$ws.Range("E34").Value = "  -7.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.68"
$ws.Range("E35").Value = "  -6.81%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0950"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "621.11"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.412"
$ws.Range("E38").Value = "  -5.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("E41").Value = "  +9.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.16"
$ws.Range("E43").Value = "  -10.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0452"
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.77"
$ws.Range("E45").Value = "  -11.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.63"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.139"
$ws.Range("E47").Value = "  -7.43%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -16.14%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.25"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.800.15"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("E51").Value = "  +0.00%  "
